# Commit: "added can be tough"
#
# 1. Appends the run " can be tough" right after the existing "s" run
#    (before the pre-existing _GoBack bookmark) in the second paragraph,
#    so "Promise" + "s" + " can be tough" reads "Promises can be tough".
# 2. Marks every paragraph (paragraph-mark run properties + every literal
#    run) as English (U.S.) - <w:lang w:val="en-US"/> - matching a
#    "Set Proofing Language -> English (United States)" pass over the
#    whole document.
#
# Range.LanguageID only ever lands on the literal runs it covers (never on
# a paragraph's own mark / w:pPr/w:rPr), so the paragraph-mark language is
# applied by rewriting each paragraph through Range.InsertXML with an
# explicit <w:rPr><w:lang .../></w:rPr> in both w:pPr and every w:r - the
# same end state Word reaches when you select the whole document and set
# its language.

$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$langRpr = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

function Escape-Xml($text) {
    $safe = $text -replace "&", "&amp;"
    $safe = $safe -replace "<", "&lt;"
    $safe = $safe -replace ">", "&gt;"
    return $safe
}

function New-Run($text) {
    $safe = Escape-Xml $text
    $space = ""
    if ($text -ne $text.Trim()) {
        $space = ' xml:space="preserve"'
    }
    return "<w:r>$langRpr<w:t$space>$safe</w:t></w:r>"
}

# --- Paragraph 1: "This is a test" -----------------------------------
$p1 = $d.Paragraphs(1)
$p1xml = "<w:p $wns><w:pPr>$langRpr</w:pPr>" + (New-Run "This is a test") + "</w:p>"
$p1.Range.InsertXML($p1xml)

# --- Paragraph 2: "Promise" + "s" + " can be tough" -------------------
$p2 = $d.Paragraphs(2)
$p2xml = "<w:p $wns><w:pPr>$langRpr</w:pPr>" `
    + (New-Run "Promise") `
    + (New-Run "s") `
    + (New-Run " can be tough") `
    + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' `
    + "</w:p>"
$p2.Range.InsertXML($p2xml)
